# Insert a new data row at row 313 (shifts existing rows 313.. down by one)
# and populate it with the new record's values, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 313, pushing old row 313 (and below) down to 314.
$ws.Rows.Item(313).Insert()

# Populate the newly inserted row 313 with the new data record.
$ws.Cells.Item(313, 1).Value  = 9
$ws.Cells.Item(313, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(313, 3).Value  = "Metropolitana"
$ws.Cells.Item(313, 4).Value  = 44809
$ws.Cells.Item(313, 5).Value  = 13
$ws.Cells.Item(313, 6).Value  = 100112032
$ws.Cells.Item(313, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(313, 8).Value  = "Sin especificar"
$ws.Cells.Item(313, 9).Value  = "Primera"
$ws.Cells.Item(313, 10).Value = 290
$ws.Cells.Item(313, 11).Value = 23000
$ws.Cells.Item(313, 12).Value = 25000
$ws.Cells.Item(313, 13).Value = 24052
$ws.Cells.Item(313, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(313, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(313, 16).Value = 481
$ws.Cells.Item(313, 17).Value = 50
$ws.Cells.Item(313, 18).Value = "Hortaliza"
